# "them chuc nang thong bao va xem thong bao"
# Adds 3-4 "notification" rows (STT, Ngay, Tieu de "xxx", Noi dung "yyy") to each
# of the gp01A / gp01B / gp02A / gp02B sheets, formats the date column, and
# sets that column's width.

$wb = $excel.ActiveWorkbook

function Fill-Notifications($SheetName, $Dates, $SelectionLastRow) {
    $ws = $wb.Worksheets.Item($SheetName)

    $row = 2
    $stt = 1
    foreach ($d in $Dates) {
        $ws.Cells.Item($row, 1).Value = $stt
        $ws.Cells.Item($row, 2).Value = $d
        $ws.Cells.Item($row, 2).NumberFormat = "d-mmm-yy"
        $ws.Cells.Item($row, 3).Value = "xxx"
        $ws.Cells.Item($row, 4).Value = "yyy"
        $row++
        $stt++
    }

    $ws.Columns.Item(2).ColumnWidth = 8.6

    $ws.Range("A2:D$SelectionLastRow").Select() | Out-Null
}

Fill-Notifications "gp01A" @(45272, 45272, 45272) 3
Fill-Notifications "gp01B" @(45272, 45273, 45274, 45275) 5
Fill-Notifications "gp02A" @(45270, 45273, 45276, 45279) 5
Fill-Notifications "gp02B" @(45270, 45271, 45272, 45273) 5

# Restore the originally active sheet/tab and its selection
$wsFinal = $wb.Worksheets.Item("gp02B")
$wsFinal.Activate() | Out-Null
$wsFinal.Range("A2:D5").Select() | Out-Null
